$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.891.62"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.585.60"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'210.38"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.476"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'0.0613"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "'18.15"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.804.48"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "1.581.58"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'0.505"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "25.880.47"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "0.0₃0724"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'60.14"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'194.12"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'4.20"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "'140.81"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").Value = "'1.70"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "'15.08"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "'6.45"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'3.03"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'0.502"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'0.778"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("D42").Value = "'0.797"
$ws.Range("E42").Value = "  +6.90%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "'5.10"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "1.718.37"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").Value = "'53.15"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -0.71%  "
